$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 12.451
$ws.Range("H3").Value = 11.495
$ws.Range("I3").Value = 0.019

$wb.Save()
